$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "members": add an "ID" column in front, renumber 1..10, tweak a
# couple of point totals, and append two new member rows.
# ---------------------------------------------------------------------------
$members = $wb.Worksheets.Item("members")

# Shift StudentID/Name/Points from A/B/C to B/C/D, insert a blank column A.
$members.Columns.Item(1).Insert()

# New header for column A.
$members.Cells.Item(1, 1).Value = "ID"
# Match the bold/centered header formatting used by the other header cells.
$members.Cells.Item(1, 2).Copy()
$members.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Sequential member IDs for the 10 existing members.
for ($r = 2; $r -le 11; $r++) {
  $members.Cells.Item($r, 1).Value = $r - 1
}

# Point total corrections.
$members.Cells.Item(2, 4).Value = 120   # Ana Reyes: 128 -> 120
$members.Cells.Item(8, 4).Value = 155   # Grace Uy: 160 -> 155

# New rows appended via the "add member" form (no sequential ID assigned yet).
$members.Cells.Item(12, 2).Value = 123145
$members.Cells.Item(12, 3).Value = "Trek II"
$members.Cells.Item(12, 4).Value = 0

$members.Cells.Item(13, 2).NumberFormat = "@"
$members.Cells.Item(13, 2).Value = "123146"
$members.Cells.Item(13, 3).Value = "Trek III"
$members.Cells.Item(13, 4).Value = 0

# ---------------------------------------------------------------------------
# Sheet "event_attendance": log attendance for a new "Mass" event.
# ---------------------------------------------------------------------------
$attendance = $wb.Worksheets.Item("event_attendance")

$massStudentIds = @(224892, 225814, 208456, 231045, 219073)
$row = 10
foreach ($studentId in $massStudentIds) {
  $attendance.Cells.Item($row, 1).Value = "Mass"
  $attendance.Cells.Item($row, 2).Value = $studentId
  $row++
}
